# Update "想去人数" (number of people interested) counts on the
# "展览" (Exhibition) and "全部类型" (All types) sheets.
#
#   F4: 1   -> 2     (丽水·thp01～风摄少微)
#   F6: 157 -> 158   (丽水·AEO纯白礼赞动漫嘉年华)

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 2
    $ws.Range("F6").Value = 158
}
